$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data (already sorted descending by total_registros), matching the
# target state of the sheet after the upload/edit.
$data = @(
    @("RUIZ CHIROQUE CLAUDIA JUDITH", 204),
    @("BANCAYAN FIESTA DILVER HUMBERTO", 179),
    @("FABIANA REBECA ARRUNATEGUI SILUPU", 176),
    @("LLENQUE ANTON HELEN JOHANA", 175),
    @("TEMOCHE ECHE URSULA YESSENIA", 171),
    @("GONZALES FIESTAS MARIA MARIBEL", 166),
    @("PINTADO CHASQUERO ESTEFANY", 166),
    @("BAUTISTA CHAVESTA ERICKA MEDALIT", 161),
    @("ANTON INGA FATIMA DEL ROSARIO", 160),
    @("VELASCO PEÑA KAREN ARELLYS", 159),
    @("HERNANDEZ CARNERO ARTURO SEBASTIAN", 152),
    @("MONDRAGON NONAJULCA MARISOL", 145),
    @("FLORES SILUPU MARY CARMEN", 145),
    @("ORDINOLA JIBAJA JOSE ALBERTO", 127),
    @("CASTRO ESTRADA CINTHIA PATRICIA", 106),
    @("MORENO YANAYACO NAYLA GUADALUPE", 104),
    @("PINTADO BENITES CRISTOBAL RODRIGO", 1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
